$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9 (Ano 2025) with refreshed faturamento data
$ws.Range("B9").Value = 3674831.7
$ws.Range("C9").Value = 578395.24
$ws.Range("D9").Value = 4253226.94
$ws.Range("E9").Value = 13.59897433547244
$ws.Range("F9").Value = 86.40102566452755
$ws.Range("G9").Value = -44.10098406562436
$ws.Range("H9").Value = -33.63761899595248
$ws.Range("I9").Value = 37177
$ws.Range("J9").Value = 1581
$ws.Range("K9").Value = 38758
$ws.Range("L9").Value = 26779
$ws.Range("M9").Value = 158.8269517158968
$ws.Range("N9").Value = 8.434394880818274
